$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Devin Booker'
$ws.Range("B2").Value = 'PG,SG'
$ws.Range("C2").Value = 'Phoenix Suns'
$ws.Range("A3").Value = 'Josh Okogie'
$ws.Range("B3").Value = 'SG,SF'
$ws.Range("C3").Value = 'Charlotte Hornets'
$ws.Range("A4").Value = 'Kawhi Leonard'
$ws.Range("B4").Value = 'SG,SF,PF'
$ws.Range("C4").Value = 'LA Clippers'
$ws.Range("A5").Value = 'Christian Braun'
$ws.Range("B5").Value = 'SG,SF'
$ws.Range("C5").Value = 'Denver Nuggets'
$ws.Range("A6").Value = 'Norman Powell'
$ws.Range("B6").Value = 'SG,SF'
$ws.Range("C6").Value = 'LA Clippers'
$ws.Range("A7").Value = 'D''Angelo Russell'
$ws.Range("B7").Value = 'PG'
$ws.Range("C7").Value = 'Brooklyn Nets'
$ws.Range("A8").Value = 'Jalen Brunson'
$ws.Range("B8").Value = 'PG'
$ws.Range("C8").Value = 'New York Knicks'
$ws.Range("A9").Value = 'Myles Turner'
$ws.Range("B9").Value = 'C'
$ws.Range("C9").Value = 'Indiana Pacers'
$ws.Range("A10").Value = 'Walker Kessler'
$ws.Range("B10").Value = 'C'
$ws.Range("C10").Value = 'Utah Jazz'
$ws.Range("A11").Value = 'Trae Young'
$ws.Range("B11").Value = 'PG'
$ws.Range("C11").Value = 'Atlanta Hawks'
$ws.Range("A12").Value = 'Desmond Bane'
$ws.Range("B12").Value = 'SG,SF'
$ws.Range("C12").Value = 'Memphis Grizzlies'
$ws.Range("A13").Value = 'Bilal Coulibaly'
$ws.Range("B13").Value = 'SG,SF'
$ws.Range("C13").Value = 'Washington Wizards'
$ws.Range("A14").Value = 'Scoot Henderson'
$ws.Range("B14").Value = 'PG'
$ws.Range("C14").Value = 'Portland Trail Blazers'
$ws.Range("A15").Value = 'LeBron James'
$ws.Range("B15").Value = 'SF,PF'
$ws.Range("C15").Value = 'Los Angeles Lakers'
$ws.Range("A16").Value = 'Jalen Williams'
$ws.Range("B16").Value = 'SG,SF,PF,C'
$ws.Range("C16").Value = 'Oklahoma City Thunder'
$ws.Range("A17").Value = 'Immanuel Quickley'
$ws.Range("B17").Value = 'PG,SG'
$ws.Range("C17").Value = 'Toronto Raptors'
$ws.Range("A18").Value = 'Brandon Ingram'
$ws.Range("B18").Value = 'SG,SF,PF'
$ws.Range("C18").Value = 'New Orleans Pelicans'
$ws.Range("A19").Value = 'Jimmy Butler'
$ws.Range("B19").Value = 'SF,PF'
$ws.Range("C19").Value = 'Miami Heat'
